{"js": "// Replace each two-digit-division problem's text with the new value.\n// Each \"before\" string is unique in the document, so a direct\n// search + replace per pair is safe and order-independent.\nconst replacements = [\n  [\"21\u00f77=3, 0\", \"26\u00f79=2, 8\"],\n  [\"66\u00f78=8, 2\", \"35\u00f76=5, 5\"],\n  [\"86\u00f79=9, 5\", \"60\u00f78=7, 4\"],\n  [\"28\u00f76=4, 4\", \"88\u00f75=17, 3\"],\n  [\"26\u00f77=3, 5\", \"51\u00f77=7, 2\"],\n  [\"29\u00f74=7, 1\", \"30\u00f78=3, 6\"],\n  [\"48\u00f73=16, 0\", \"14\u00f78=1, 6\"],\n  [\"36\u00f77=5, 1\", \"74\u00f77=10, 4\"],\n  [\"16\u00f74=4, 0\", \"50\u00f76=8, 2\"],\n  [\"44\u00f74=11, 0\", \"24\u00f75=4, 4\"],\n  [\"11\u00f76=1, 5\", \"18\u00f73=6, 0\"],\n  [\"66\u00f72=33, 0\", \"64\u00f79=7, 1\"],\n  [\"41\u00f73=13, 2\", \"92\u00f75=18, 2\"],\n  [\"27\u00f75=5, 2\", \"76\u00f72=38, 0\"],\n  [\"45\u00f72=22, 1\", \"31\u00f75=6, 1\"],\n  [\"89\u00f77=12, 5\", \"10\u00f76=1, 4\"],\n  [\"16\u00f75=3, 1\", \"10\u00f74=2, 2\"],\n  [\"41\u00f75=8, 1\", \"52\u00f74=13, 0\"],\n  [\"58\u00f78=7, 2\", \"96\u00f74=24, 0\"],\n  [\"52\u00f77=7, 3\", \"57\u00f77=8, 1\"],\n  [\"76\u00f76=12, 4\", \"69\u00f73=23, 0\"],\n  [\"72\u00f77=10, 2\", \"25\u00f76=4, 1\"],\n  [\"53\u00f79=5, 8\", \"42\u00f76=7, 0\"],\n  [\"93\u00f78=11, 5\", \"62\u00f72=31, 0\"],\n  [\"88\u00f73=29, 1\", \"88\u00f72=44, 0\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n", "ps1": "# Update the twenty-five two-digit-divided-by-one-digit problems in the\n# table with their new dividend/divisor/quotient/remainder values.\n# Every \"old\" string below is unique in the document body, so a plain\n# Find/Replace (wdReplaceOne) per pair is safe and order-independent.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"21\u00f77=3, 0\", \"26\u00f79=2, 8\"),\n    @(\"66\u00f78=8, 2\", \"35\u00f76=5, 5\"),\n    @(\"86\u00f79=9, 5\", \"60\u00f78=7, 4\"),\n    @(\"28\u00f76=4, 4\", \"88\u00f75=17, 3\"),\n    @(\"26\u00f77=3, 5\", \"51\u00f77=7, 2\"),\n    @(\"29\u00f74=7, 1\", \"30\u00f78=3, 6\"),\n    @(\"48\u00f73=16, 0\", \"14\u00f78=1, 6\"),\n    @(\"36\u00f77=5, 1\", \"74\u00f77=10, 4\"),\n    @(\"16\u00f74=4, 0\", \"50\u00f76=8, 2\"),\n    @(\"44\u00f74=11, 0\", \"24\u00f75=4, 4\"),\n    @(\"11\u00f76=1, 5\", \"18\u00f73=6, 0\"),\n    @(\"66\u00f72=33, 0\", \"64\u00f79=7, 1\"),\n    @(\"41\u00f73=13, 2\", \"92\u00f75=18, 2\"),\n    @(\"27\u00f75=5, 2\", \"76\u00f72=38, 0\"),\n    @(\"45\u00f72=22, 1\", \"31\u00f75=6, 1\"),\n    @(\"89\u00f77=12, 5\", \"10\u00f76=1, 4\"),\n    @(\"16\u00f75=3, 1\", \"10\u00f74=2, 2\"),\n    @(\"41\u00f75=8, 1\", \"52\u00f74=13, 0\"),\n    @(\"58\u00f78=7, 2\", \"96\u00f74=24, 0\"),\n    @(\"52\u00f77=7, 3\", \"57\u00f77=8, 1\"),\n    @(\"76\u00f76=12, 4\", \"69\u00f73=23, 0\"),\n    @(\"72\u00f77=10, 2\", \"25\u00f76=4, 1\"),\n    @(\"53\u00f79=5, 8\", \"42\u00f76=7, 0\"),\n    @(\"93\u00f78=11, 5\", \"62\u00f72=31, 0\"),\n    @(\"88\u00f73=29, 1\", \"88\u00f72=44, 0\"),\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $found = $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n    if (-not $found) {\n        throw \"Text not found: $old\"\n    }\n}\n"}
